$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = "<kilo>"
$ws.Range("C2").Value = 23

# Row 3
$ws.Range("C3").Value = 28

# Row 4
$ws.Range("C4").Value = 31

# Row 5
$ws.Range("B5").Value = "<line>"
$ws.Range("C5").Value = 37

# Row 6
$ws.Range("C6").Value = 36

# Row 7
$ws.Range("C7").Value = 35

# Row 8
$ws.Range("C8").Value = 36

# Row 9
$ws.Range("B9").Value = "<nove>"
$ws.Range("C9").Value = 25

# Row 10
$ws.Range("C10").Value = 30

# Row 11
$ws.Range("C11").Value = 35

# Row 12
$ws.Range("C12").Value = 38

# Row 13
$ws.Range("C13").Value = 27

# Row 14
$ws.Range("C14").Value = 32

# Row 15
$ws.Range("B15").Value = "<can>"
$ws.Range("C15").Value = 14
